# Updated cryptos list (price / 1h volume change refresh + TRON/WrappedEther row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '69.295.73'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '3.384.87'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'" + '588.12'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = "'" + '180.60'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'" + '0.594'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').Value = "'" + '0.195'
$ws.Range('E9').Value = '  +8.56%  '
$ws.Range('D10').Value = "'" + '0.588'
$ws.Range('D11').Value = "'" + '48.80'
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('D12').Value = "'" + '0.0000285'
$ws.Range('E12').Value = '  +5.03%  '
$ws.Range('D13').Value = "'" + '683.49'
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('D14').Value = "'" + '8.64'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '3.930.11'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '69.324.42'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.385.45'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'" + '0.120'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').Value = "'" + '17.78'
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('D20').Value = "'" + '11.41'
$ws.Range('E20').Value = '  +3.87%  '
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = "'" + '17.06'
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('D24').Value = "'" + '104.53'
$ws.Range('E24').Value = '  +6.35%  '
$ws.Range('D25').Value = "'" + '3.94'
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').Value = "'" + '9.61'
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').Value = "'" + '34.51'
$ws.Range('E28').Value = '  +3.79%  '
$ws.Range('D29').Value = "'" + '8.69'
$ws.Range('E29').Value = '  +1.84%  '
$ws.Range('D30').Value = "'" + '7.00'
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').Value = "'" + '11.20'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').Value = "'" + '558.68'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('E33').Value = '  +9.28%  '
$ws.Range('D34').Value = "'" + '0.106'
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('D35').Value = "'" + '58.15'
$ws.Range('E35').Value = '  +1.52%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = '3.714.81'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  +8.18%  '
$ws.Range('D39').Value = "'" + '34.96'
$ws.Range('E39').Value = '  +2.74%  '
$ws.Range('D40').Value = "'" + '3.26'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('D41').Value = '0.0₃0706'
$ws.Range('E41').Value = '  +4.62%  '
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').Value = "'" + '0.340'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').Value = "'" + '0.0418'
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('D45').Value = "'" + '3.27'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').Value = "'" + '2.66'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('E48').Value = '  +5.90%  '
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('D50').Value = "'" + '132.57'
$ws.Range('E50').Value = '  +3.11%  '
$ws.Range('E51').Value = '  -3.18%  '
